$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: fill in the Bollinger / PriceChange / UpDown columns that were
# left blank before (the predictor now also scores this day) ---
$ws.Range("W3").Value = 0
$ws.Range("X3").Value = 0.19000099999999875
$ws.Range("Y3").Value = "Up"

# --- Row 4: brand new day of predictions appended to the table ---
$ws.Range("A4").Value = 42641.890729166669
$ws.Range("B4").Value = -2
$ws.Range("C4").Value = "Neutral"
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = "Random"
$ws.Range("Q4").Value = 65.63785237683328
$ws.Range("R4").Value = 0.48
$ws.Range("S4").Value = 0.068500000000000005
$ws.Range("T4").Value = -0.061899999999999997
$ws.Range("U4").Value = 2.27
$ws.Range("V4").Value = "N/A"
$ws.Range("W4").Value = 0

# Match formatting of the sibling rows above instead of inventing new
# number formats: A4 should carry the same date style as A2/A3, and
# S4/T4 the same percentage style as S2:T3.
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("S3:T3").Copy()
$ws.Range("S4:T4").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$wb.Save()
